$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.648.81"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "3.940.80"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.60"
$ws.Range("E5").Value = "  +8.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.00"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.729"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.85"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.44"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").Value = "4.568.62"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "3.943.11"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.07"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  +7.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.91"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "69.585.48"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.59"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.59"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  +13.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.95"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.74"
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "702.20"
$ws.Range("E29").Value = "  -3.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.36"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "68.46"
$ws.Range("E33").Value = "  +12.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.454"
$ws.Range("E34").Value = "  +12.71%  "
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "0.0₃0871"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.59"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.150"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("E43").Value = "  +5.92%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.01"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.24"
$ws.Range("E45").Value = "  +15.26%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.143"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "0.0₆0360"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.34"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.14"
$ws.Range("E51").Value = "  +0.17%  "
